$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheets 2..6 (tab order "Sheet4", "Sheet3", "Sheet2", "Sheet1") each get a
# new header row inserted at the top: country / search index / <keyword>.
# These are processed first (in an order that does NOT matter for the final
# "active sheet", since Sheet5 -- processed last below -- becomes the
# selected tab).
# ---------------------------------------------------------------------------

# Worksheets.Item(3) = "Sheet4" -> xl/worksheets/sheet3.xml (keyword: generative ai)
$ws = $wb.Worksheets.Item(3)
$ws.Rows.Item(1).Insert()
$ws.Range("A1").Value = "country"
$ws.Range("B1").Value = "search index"
$ws.Range("C1").Value = "generative ai"
$ws.Rows.Item(1).Select()

# Worksheets.Item(4) = "Sheet3" -> xl/worksheets/sheet4.xml (keyword: ML)
$ws = $wb.Worksheets.Item(4)
$ws.Rows.Item(1).Insert()
$ws.Range("A1").Value = "country"
$ws.Range("B1").Value = "search index"
$ws.Range("C1").Value = "ML"
$ws.Rows.Item(1).Select()

# Worksheets.Item(5) = "Sheet2" -> xl/worksheets/sheet5.xml (keyword: chatgpt)
$ws = $wb.Worksheets.Item(5)
$ws.Rows.Item(1).Insert()
$ws.Range("A1").Value = "country"
$ws.Range("B1").Value = "search index"
$ws.Range("C1").Value = "chatgpt"
$ws.Rows.Item(1).Select()

# Worksheets.Item(6) = "Sheet1" -> xl/worksheets/sheet6.xml (keyword: deep learning)
$ws = $wb.Worksheets.Item(6)
$ws.Rows.Item(1).Insert()
$ws.Range("A1").Value = "country"
$ws.Range("B1").Value = "search index"
$ws.Range("C1").Value = "deep learning"
$ws.Rows.Item(1).Select()

# Worksheets.Item(2) = "Sheet5" -> xl/worksheets/sheet2.xml (keyword: datascience)
# Processed LAST so that it ends up as the active / selected tab, matching
# the target workbook (bookViews activeTab="1", sheetView tabSelected="1").
$ws = $wb.Worksheets.Item(2)
$ws.Rows.Item(1).Insert()
$ws.Range("A1").Value = "country"
$ws.Range("B1").Value = "search index"
$ws.Range("C1").Value = "datascience"
$ws.Range("J8").Select()
